# Commit: Tue, May 05, 2020  8:06:54 PM
#
# 1) The table on the "SOURCES OF FINANCE" slide had its table style
#    switched from the deck's custom style ({EA188C00-...}, the
#    tableStyles.xml default) to the built-in "Medium Style 2 - Accent 1"
#    style ({289E0716-AD4E-4899-9AFE-E035ACB3935A}).
#
# 2) The Slide Master / Notes Master theme assignment was swapped (the
#    deck's "Integral" theme and the default "Office Theme" traded
#    places between the two masters). We call the closest PowerPoint
#    object-model equivalents for that (Master.ApplyTheme) for fidelity,
#    even though the underlying theme-part content is managed by the
#    host application rather than by script-visible state.

$p = $ppt.ActivePresentation

# --- 1. Re-style the table on the slide that lists sources of finance ---
$styleApplied = $false
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle("{289E0716-AD4E-4899-9AFE-E035ACB3935A}")
            $styleApplied = $true
        }
    }
}

# --- 2. Swap the theme carried by the slide master and the notes master ---
$master = $p.SlideMaster
$notesMaster = $p.NotesMaster

$masterTheme = $master.Theme
$notesMasterTheme = $notesMaster.Theme

$master.ApplyTheme($notesMasterTheme)
$notesMaster.ApplyTheme($masterTheme)

Write-Output "table style applied: $styleApplied"
